# Applies the "client_portfolio_report" refresh: updated generation
# timestamp, updated trade totals/summary, and refreshed "Recent Trades"
# listing (new trades, reordered/retimed entries).

$d = $word.ActiveDocument

function Replace-ParagraphText($Index, $OldText, $NewText) {
    $p = $d.Paragraphs($Index)
    $r = $p.Range
    $ok = $r.Find.Execute($OldText, $true, $false, $false, $false, $false,
                           $true, 1, $false, $NewText, 2)
    if (-not $ok) {
        throw "Replace failed for paragraph $Index (old text not found): $OldText"
    }
}

# Header / generated timestamp
Replace-ParagraphText 2 `
    "Generated: 2025-11-23 05:03 PM" `
    "Generated: 2025-11-23 10:04 PM"

# Trade Summary
Replace-ParagraphText 5 `
    "Total Trades: 27" `
    "Total Trades: 41"

# Portfolio Analysis narrative
Replace-ParagraphText 7 `
    "Portfolio contains 27 total trades (8 buys, 15 sells) with total transaction value of `$1,297,500.00. This report provides a comprehensive overview of client trading activity for compliance review and regulatory monitoring." `
    "Portfolio contains 41 total trades (25 buys, 16 sells) with total transaction value of `$6,626,917.50. This report provides a comprehensive overview of client trading activity for compliance review and regulatory monitoring."

# Recent Trades listing (each line replaced in-place by paragraph index so
# duplicate source lines map to their correct, distinct replacements)
Replace-ParagraphText 9 `
    "2025-11-23 01:35 PM - Sheila Carter: Sell 500 TSLA @ `$180.0" `
    "2025-11-23 05:30:00 - Sandra Carter: Sell 225 EBAY @ `$295.5"

Replace-ParagraphText 10 `
    "2025-11-23 01:35 PM - Sheila Carter: Sell 500 TSLA @ `$180" `
    "2025-11-23 08:00:00 - Timothy Mitchell: Buy 270 SPOT @ `$0.0"

Replace-ParagraphText 11 `
    "2025-11-23 03:13 PM - Tony Stark: Buy 2500 PLTR @ `$25.5" `
    "2025-11-23 10:30:00 - Catherine Perez: Buy 340 ROKU @ `$142.25"

Replace-ParagraphText 12 `
    "2025-11-23 03:13 PM - Tony Stark: Buy 2500 PLTR @ `$25.5" `
    "2025-11-23 13:00:00 - Gregory Roberts: Sell 195 PINS @ `$0.0"

Replace-ParagraphText 13 `
    "2025-11-23 03:54 PM - Sheila Carter: Sell 500 TSLA @ `$180.0" `
    "2025-11-23 09:44 PM - Bruce Wayne: Buy 10000 GOOGL @ `$142.5"

Replace-ParagraphText 14 `
    "2025-11-23 03:55 PM - Sheila Carter: Sell 500 TSLA @ `$180.0" `
    "2025-11-23 09:44 PM - Bruce Wayne: Buy 10000 GOOGL @ `$142.5"

Replace-ParagraphText 15 `
    "2025-11-23 03:57 PM - Sheila Carter: Sell 500 TSLA @ `$180.0" `
    "2025-11-23 09:47 PM - Sheila Carter: Sell 500 TSLA @ `$180.0"

Replace-ParagraphText 16 `
    "2025-11-23 04:08 PM - Mr Jones: Buy 1 THREE SHOTS FOUR ENDOSYMBIOSIS @ `$0" `
    "2025-11-23 09:47 PM - Sheila Carter: SELL 500 TSLA @ `$180"

Replace-ParagraphText 17 `
    "2025-11-23 04:57 PM - Sheila Carter: Sell 500 TSLA @ `$180.0" `
    "2025-11-23 10:00 PM - Bruce Wayne: Buy 10000 GOOGL @ `$142.5"

Replace-ParagraphText 18 `
    "2025-11-23 04:57 PM - Sheila Carter: SELL 500 TSLA @ `$180" `
    "2025-11-23 10:00 PM - Bruce Wayne: Buy 10000 GOOGL @ `$142.5"
